$wb = $excel.ActiveWorkbook

# "Согласие" sheet - concordance matrix
$wsConcordance = $wb.Worksheets.Item("Согласие")
$wsConcordance.Range("C2").Value = 0.6875
$wsConcordance.Range("D2").Value = 0.9375
$wsConcordance.Range("B3").Value = 0.3125
$wsConcordance.Range("D3").Value = 0.9375
$wsConcordance.Range("B4").Value = 0.0625
$wsConcordance.Range("C4").Value = 0.0625

# "Несогласие" sheet - discordance matrix
$wsDiscordance = $wb.Worksheets.Item("Несогласие")
$wsDiscordance.Range("C2").Value = 0.0555555555555556
$wsDiscordance.Range("D2").Value = 0.0666666666666667
$wsDiscordance.Range("D3").Value = 0.1
$wsDiscordance.Range("C4").Value = 0.5
